$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: simple text replacements on the first three rows ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- Step 2: insert 10 new rows right after row 3 (i.e. before the old row 4) ---
# Rows.Add(beforeRow) always inserts immediately above the fixed anchor, so
# walking the desired values in reverse yields the correct final top-to-bottom order.
$newValues = @("102", "0.00003", "0.00008", "0.00005", "0.00001", "0.00005", "0.00005", "0.00005", "0.00402", "100.0")

$anchorRow = $t.Rows.Item(4)
$count = $newValues.Count
for ($i = $count - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($anchorRow)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
}

# --- Step 3: clean up the two multi-run (tab separated) data rows ---
# Old row 34 (1-indexed) -> now shifted down by 10 new rows = row 44
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"

# Old row 35 (1-indexed) -> now row 45
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0"

# Old row 36 (1-indexed, previously empty) -> now row 46
$t.Rows.Item(46).Cells.Item(1).Range.Text = "105"
